$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
$ws.Range("L1").Value = "Travis R."
